# Loan RBI, Variable Instalments
# The "Repayment schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column so a variable-instalment column can be
# tracked separately from the "Outstanding" figures further down the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)   # "Repayment schedule"

# Insert a new blank column at N (pushes Late / heading / Outstanding one
# column to the right: N->O, O->P, P->Q) and give it the same width as the
# column immediately to its left ("M"), matching what Excel does when you
# insert a column from the sheet grid.
$leftWidth = $ws.Range("M1").ColumnWidth
$ws.Range("N1").EntireColumn.Insert()
$ws.Range("N1").ColumnWidth = $leftWidth

# Bring the "Repayment schedule" tab to the front and leave the selection
# where the author left it while reviewing the new column.
$ws.Activate()
[void]$ws.Range("S8").Select()
